$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Status and Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B6").Value = "active"
$wsMeta.Range("B8").Value = "2026-01-28T10:29:57+00:00"

# --- Sheet "Concepts": replace the concept list with the new GT Structure codes ---
$wsConcepts = $wb.Worksheets.Item("Concepts")

$codes = @(
    @("BAT", "Bâtiment"),
    @("ETAG", "Étage"),
    @("COUL", "Couloir"),
    @("AILE", "Aile"),
    @("BOX", "Box"),
    @("CHAMB", "Chambre"),
    @("LIT", "lit"),
    @("PL_TECH", "Plateau technique"),
    @("PNT_CLCT", "Point de collecte"),
    @("PNT_LVRSN", "Point de livraison"),
    @("SL_EXM", "Salle examen"),
    @("SL_CONS", "Salle de consultation")
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = $i + 2
    $code = $codes[$i][0]
    $display = $codes[$i][1]

    $wsConcepts.Cells.Item($row, 2).Value = $code
    $wsConcepts.Cells.Item($row, 3).Value = $display
    $wsConcepts.Cells.Item($row, 4).Value = ""
}
